$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the weekly report row for 07/01/2025 (D48:J48)
$ws.Range("D48").Value = (Get-Date -Year 2025 -Month 1 -Day 7).Date
$ws.Range("E48").Value = 192
$ws.Range("F48").Value = 734
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 1012
$ws.Range("J48").Value = "N/A"

# Update the view state to reflect where the author left off
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 5
$ws.Range("J49").Select()
